# "run all scenario sulawesi selatan"
#
# This edit:
#  1. Inserts a new column D on Sheet1 that applies a +0.3% adjustment
#     (=C*1.003) to the LRC rate computed in column C, with a handful of
#     rows overridden by hand-typed literal values (matching the source
#     edit), shifting the old D:F (Jabar sector name / his rate / rate2)
#     columns to E:G.
#  2. Appends a new row (62) for "sektor lainnya" with zeroed rate values.
#  3. Adds the new shared string "sektor lainnya".
#  4. Updates the selection / active sheet so Sheet1 (not Sheet3) is now
#     the active tab, matching the "run all scenario" pass over the data.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------
# 1. Insert new column D (pushes old D/E/F -> E/F/G) and fill it with the
#    1.003 uplift formula for every data row (2-61).
# ---------------------------------------------------------------------
$ws1.Columns.Item(4).Insert()

$ws1.Range("D2:D61").Formula = "=C2*1.003"

# A handful of rows in the source workbook were overwritten by hand with
# literal numbers instead of the dragged-down formula -- reproduce those
# exact overrides.
$ws1.Range("D3").Value = 1
$ws1.Range("D6").Value = 1
$ws1.Range("D24").Value = 0.81
$ws1.Range("D25").Value = 0.81
$ws1.Range("D26").Value = 0.81
$ws1.Range("D27").Value = 0.81
$ws1.Range("D28").Value = 0.81
$ws1.Range("D29").Value = 0.81
$ws1.Range("D30").Value = 0.81
$ws1.Range("D31").Value = 1
$ws1.Range("D59").Value = 0.8

# ---------------------------------------------------------------------
# 2 & 3. New row 62: "sektor lainnya" with zero rates.
# ---------------------------------------------------------------------
$ws1.Range("B62").Value = "sektor lainnya"
$ws1.Range("C62").Value = 0
$ws1.Range("D62").Value = 0

# ---------------------------------------------------------------------
# 4. Selection / active sheet bookkeeping.
# ---------------------------------------------------------------------
$ws3.Range("C1:C1048576").Select()
$ws3.Range("C31").Activate()

$ws1.Activate()
$ws1.Range("C2:C62").Select()
$ws1.Range("C2").Activate()
